# Auto-generated Excel COM-interop script
# Applies market-data refresh values (currentAveragePrice* / LevePrice* / LeveProfit*)
# to the leve-profit tables on each class sheet, per the scheduled-runner update.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 26317428
$ws.Range("I33").Value = 35716052
$ws.Range("K33").Value = 35716052
$ws.Range("M33").Value = -35715823
$ws.Range("H137").Value = 3687.2354
$ws.Range("I137").Value = 3290
$ws.Range("K137").Value = 9870
$ws.Range("M137").Value = -7320
$ws.Range("H138").Value = 1727789.4
$ws.Range("I138").Value = 640.67645
$ws.Range("J138").Value = 4174583.2
$ws.Range("K138").Value = 1922.02935
$ws.Range("L138").Value = 12523749.6
$ws.Range("M138").Value = 3217.97065
$ws.Range("N138").Value = -12534029.6
$ws.Range("H141").Value = 1359.2916
$ws.Range("I141").Value = 1201.174
$ws.Range("K141").Value = 3603.522
$ws.Range("M141").Value = 1576.478

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 55558910
$ws.Range("I2").Value = 1285.091
$ws.Range("K2").Value = 1285.091
$ws.Range("M2").Value = -1172.091
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -701
$ws.Range("H45").Value = 3939.3125
$ws.Range("I45").Value = 1478
$ws.Range("K45").Value = 1478
$ws.Range("M45").Value = -1101
$ws.Range("H61").Value = 5454.951
$ws.Range("I61").Value = 2024.069
$ws.Range("K61").Value = 2024.069
$ws.Range("M61").Value = -1812.069
$ws.Range("H74").Value = 24738.25
$ws.Range("I74").Value = 31374.809
$ws.Range("J74").Value = 4828.5713
$ws.Range("K74").Value = 31374.809
$ws.Range("L74").Value = 4828.5713
$ws.Range("M74").Value = -30500.809
$ws.Range("N74").Value = -6576.5713
$ws.Range("H77").Value = 24738.25
$ws.Range("I77").Value = 31374.809
$ws.Range("J77").Value = 4828.5713
$ws.Range("K77").Value = 156874.045
$ws.Range("L77").Value = 24142.8565
$ws.Range("M77").Value = -152506.045
$ws.Range("N77").Value = -32878.85649999999
$ws.Range("H95").Value = 40266.75
$ws.Range("J95").Value = 40266.75
$ws.Range("L95").Value = 40266.75
$ws.Range("N95").Value = -45758.75
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("H113").Value = 56546
$ws.Range("J113").Value = 56546
$ws.Range("L113").Value = 56546
$ws.Range("N113").Value = -65224
$ws.Range("H116").Value = 55558910
$ws.Range("I116").Value = 1285.091
$ws.Range("K116").Value = 1285.091
$ws.Range("M116").Value = 1008.909
$ws.Range("H132").Value = 5316.1665
$ws.Range("I132").Value = 2413.1482
$ws.Range("K132").Value = 7239.444600000001
$ws.Range("M132").Value = -4709.444600000001
$ws.Range("H136").Value = 5454.951
$ws.Range("I136").Value = 2024.069
$ws.Range("K136").Value = 6072.207
$ws.Range("M136").Value = -3522.207
$ws.Range("N96").ClearContents()

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 55558910
$ws.Range("I3").Value = 1285.091
$ws.Range("K3").Value = 1285.091
$ws.Range("M3").Value = -1171.091
$ws.Range("H105").Value = 2917.1667
$ws.Range("I105").Value = 2154.353
$ws.Range("K105").Value = 2154.353
$ws.Range("M105").Value = -407.3530000000001
$ws.Range("H107").Value = 80362460
$ws.Range("I107").Value = 140627490
$ws.Range("J107").Value = 9107
$ws.Range("K107").Value = 140627490
$ws.Range("L107").Value = 9107
$ws.Range("M107").Value = -140625570
$ws.Range("N107").Value = -12947
$ws.Range("H129").Value = 56269.25
$ws.Range("J129").Value = 58359.668
$ws.Range("L129").Value = 58359.668
$ws.Range("N129").Value = -68359.66800000001
$ws.Range("H134").Value = 5013.7856
$ws.Range("I134").Value = 1867.3784
$ws.Range("K134").Value = 5602.135200000001
$ws.Range("M134").Value = -3067.135200000001
$ws.Range("H140").Value = 90666.336
$ws.Range("J140").Value = 90666.336
$ws.Range("L140").Value = 90666.336
$ws.Range("N140").Value = -101026.336

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7760052
$ws.Range("J31").Value = 15163856
$ws.Range("L31").Value = 15163856
$ws.Range("N31").Value = -15164446
$ws.Range("H34").Value = 7760052
$ws.Range("J34").Value = 15163856
$ws.Range("L34").Value = 15163856
$ws.Range("N34").Value = -15164260
$ws.Range("H58").Value = 10005992
$ws.Range("I58").Value = 20836054
$ws.Range("K58").Value = 20836054
$ws.Range("M58").Value = -20835851
$ws.Range("H134").Value = 5367.7817
$ws.Range("I134").Value = 2396.889
$ws.Range("J134").Value = 8232.571
$ws.Range("K134").Value = 7190.667
$ws.Range("L134").Value = 24697.713
$ws.Range("M134").Value = -4655.667
$ws.Range("N134").Value = -29767.713
$ws.Range("H136").Value = 10005992
$ws.Range("I136").Value = 20836054
$ws.Range("K136").Value = 62508162
$ws.Range("M136").Value = -62505612

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4252.9414
$ws.Range("I5").Value = 2400.1
$ws.Range("J5").Value = 6899.857
$ws.Range("K5").Value = 7200.299999999999
$ws.Range("L5").Value = 20699.571
$ws.Range("M5").Value = -7088.299999999999
$ws.Range("N5").Value = -20923.571
$ws.Range("H8").Value = 3374.5
$ws.Range("I8").Value = 3374.5
$ws.Range("K8").Value = 10123.5
$ws.Range("M8").Value = -9984.5
$ws.Range("H39").Value = 15021.556
$ws.Range("J39").Value = 18599.143
$ws.Range("L39").Value = 55797.429
$ws.Range("N39").Value = -56385.429
$ws.Range("H68").Value = 3486.6428
$ws.Range("I68").Value = 1941.5714
$ws.Range("J68").Value = 5031.7144
$ws.Range("K68").Value = 5824.7142
$ws.Range("L68").Value = 15095.1432
$ws.Range("M68").Value = -5013.7142
$ws.Range("N68").Value = -16717.1432
$ws.Range("H71").Value = 3486.6428
$ws.Range("I71").Value = 1941.5714
$ws.Range("J71").Value = 5031.7144
$ws.Range("K71").Value = 17474.1426
$ws.Range("L71").Value = 45285.4296
$ws.Range("M71").Value = -13418.1426
$ws.Range("N71").Value = -53397.4296
$ws.Range("H92").Value = 4809139
$ws.Range("I92").Value = 886.3333
$ws.Range("J92").Value = 5918735.5
$ws.Range("K92").Value = 2658.9999
$ws.Range("L92").Value = 17756206.5
$ws.Range("M92").Value = -1410.9999
$ws.Range("N92").Value = -17758702.5
$ws.Range("H129").Value = 13889447
$ws.Range("I129").Value = 298.66666
$ws.Range("J129").Value = 55556892
$ws.Range("K129").Value = 895.9999799999999
$ws.Range("L129").Value = 166670676
$ws.Range("M129").Value = 4104.00002
$ws.Range("N129").Value = -166680676
$ws.Range("H135").Value = 4252.9414
$ws.Range("I135").Value = 2400.1
$ws.Range("J135").Value = 6899.857
$ws.Range("K135").Value = 21600.9
$ws.Range("L135").Value = 62098.713
$ws.Range("M135").Value = -19065.9
$ws.Range("N135").Value = -67168.713
$ws.Range("H140").Value = 2377.64
$ws.Range("I140").Value = 1025.0555
$ws.Range("K140").Value = 3075.1665
$ws.Range("M140").Value = 2104.8335
$ws.Range("H141").Value = 6018
$ws.Range("I141").Value = 2689.6667
$ws.Range("K141").Value = 8069.000100000001
$ws.Range("M141").Value = -2889.000100000001

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 28949.25
$ws.Range("J35").Value = 28949.25
$ws.Range("L35").Value = 28949.25
$ws.Range("N35").Value = -29545.25
$ws.Range("H44").Value = 11999.8
$ws.Range("I44").Value = 9999.75
$ws.Range("J44").Value = 20000
$ws.Range("K44").Value = 9999.75
$ws.Range("L44").Value = 20000
$ws.Range("M44").Value = -9403.75
$ws.Range("N44").Value = -21192
$ws.Range("H102").Value = 1523.1482
$ws.Range("I102").Value = 1323.6595
$ws.Range("K102").Value = 1323.6595
$ws.Range("M102").Value = 298.3405
$ws.Range("H117").Value = 46204.2
$ws.Range("J117").Value = 46204.2
$ws.Range("L117").Value = 46204.2
$ws.Range("N117").Value = -53088.2
$ws.Range("H122").Value = 6533428.5
$ws.Range("I122").Value = 14367193
$ws.Range("K122").Value = 43101579
$ws.Range("M122").Value = -43099129
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("H132").Value = 4882.1113
$ws.Range("I132").Value = 2399.9092
$ws.Range("K132").Value = 7199.7276
$ws.Range("M132").Value = -4669.7276
$ws.Range("H135").Value = 50241.11
$ws.Range("J135").Value = 50241.11
$ws.Range("L135").Value = 50241.11
$ws.Range("N135").Value = -60381.11
$ws.Range("N123").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6084.093
$ws.Range("I136").Value = 1855.1875
$ws.Range("K136").Value = 5565.5625
$ws.Range("M136").Value = -3015.5625

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 12546.091
$ws.Range("I45").Value = 999
$ws.Range("J45").Value = 13700.8
$ws.Range("K45").Value = 999
$ws.Range("L45").Value = 13700.8
$ws.Range("M45").Value = -508
$ws.Range("N45").Value = -14682.8
$ws.Range("H113").Value = 9300.379000000001
$ws.Range("I113").Value = 12052.909
$ws.Range("K113").Value = 36158.727
$ws.Range("M113").Value = -33988.727
$ws.Range("H132").Value = 13900639
$ws.Range("I132").Value = 17246098
$ws.Range("K132").Value = 51738294
$ws.Range("M132").Value = -51735764

